$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 new rows before the existing data (which starts at row 2),
# pushing the original 7 data rows down to rows 12-18.
$ws.Rows("2:11").Insert()

# The inserted rows inherit formatting from the row above (the bold
# header row), so clear that back out to match the unformatted data
# cells used elsewhere in the sheet.
$ws.Range("A2:C11").ClearFormats()

# Populate the newly inserted rows with the new accelerometer samples.
$newData = @(
  @(-0.3275403976440429, -1.032773733139038, -0.2356588244438171),
  @(-1.590312331914902, -0.9565697163343428, -1.569369990378618),
  @(0.2956193089485168, -0.9267413020133971, -1.241599485278129),
  @(-0.6574213504791258, -0.9026327282190323, -1.075559064745903),
  @(0.3167376518249511, -0.9883218407630923, -1.371818482875824),
  @(0.4182748794555663, -0.9399109184741974, -1.527341216802597),
  @(0.2231501340866089, -0.9057424068450929, -1.459591150283814),
  @(-0.4751685261726377, -0.88183431327343, -1.828994989395142),
  @(-0.972740650177002, -1.058116793632508, -2.322797894477844),
  @(-1.244342982769013, -1.110738858580589, -2.960273459553719)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 1).Value = $newData[$i][0]
  $ws.Cells.Item($row, 2).Value = $newData[$i][1]
  $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}
